# Update VAT invoice import template: add a new "Payment Status" column (E)
# with header text, matching data-cell formatting, and a comment on the
# header cell explaining the expected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (style + border) of the existing last column (D1:D2)
# into the new column E1:E2, the same way the author would have dragged the
# "Voucher ID *" header formatting rightward before typing the new header.
$ws.Range("D1:D2").Copy()
$ws.Range("E1:E2").PasteSpecial(-4122)

# New header text for the added column.
$ws.Range("E1").Value = "Payment Status"

# Size column E to fit its header text (closest the host lets us get to the
# real "best fit" pixel width).
$ws.Columns.Item(5).ColumnWidth = 14

# Document the expected values for the new column via a cell comment.
$comment = $ws.Range("E1").AddComment("Alex Phuong:`nTrạng thái thanh toán: Paid/Unpaid`n")

# Restore the selection to match the saved state of the workbook.
$ws.Range("G4").Select() | Out-Null
